# Actualización automática 2025-06-16 13:01:14
# Adds a new "GRANITO" column (inserted before the "GRIFERIAS" column) and
# three new trailing columns ("NO RESURTIBLES", "PANELES PVC", "PANELES PU")
# to the "VENTAS POR GRUPO" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# ---------------------------------------------------------------------
# 1) Insert a new column F ("GRANITO") before the existing "GRIFERIAS"
#    column. This shifts the old F:N columns one place to the right
#    (to G:O), carrying their formatting along automatically.
# ---------------------------------------------------------------------
$ws.Columns("F").Insert()

$ws.Cells.Item(1, 6).Value = "GRANITO"
$ws.Range("F2:F52").Value = 0
$ws.Cells.Item(53, 6).Value = "0 de 51"

# Excel's ColumnWidth property is offset by ~0.83 from the raw OOXML
# "width" attribute for this workbook's default font, so subtract that
# to land exactly on the target stored width.
$ws.Columns(6).ColumnWidth = 13 - 0.83

# ---------------------------------------------------------------------
# 2) Append three new trailing columns P:R ("NO RESURTIBLES",
#    "PANELES PVC", "PANELES PU") after the last existing column (now O).
#    Copy formatting from column O first so headers/data/summary rows
#    pick up the right styles, then write the actual values.
# ---------------------------------------------------------------------
$ws.Range("O1").Copy()
$ws.Range("P1:R1").PasteSpecial(-4122)

$ws.Range("O2").Copy()
$ws.Range("P2:R52").PasteSpecial(-4122)

$ws.Range("O53").Copy()
$ws.Range("P53:R53").PasteSpecial(-4122)

$ws.Cells.Item(1, 16).Value = "NO RESURTIBLES"
$ws.Cells.Item(1, 17).Value = "PANELES PVC"
$ws.Cells.Item(1, 18).Value = "PANELES PU"

$ws.Range("P2:R52").Value = 0
$ws.Cells.Item(5, 16).Value = 2.12
$ws.Cells.Item(21, 16).Value = 129.25

$ws.Cells.Item(53, 16).Value = "2 de 51"
$ws.Cells.Item(53, 17).Value = "0 de 51"
$ws.Cells.Item(53, 18).Value = "0 de 51"

$ws.Columns(16).ColumnWidth = 20 - 0.83
$ws.Columns(17).ColumnWidth = 17 - 0.83
$ws.Columns(18).ColumnWidth = 16 - 0.83
